$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standardise "cost_variable" -> "cost_variable_om" for the whole parameter
# column range that used it (C10:C39).
$ws.Range("C10:C39").Value2 = "cost_variable_om"

# Match the author's recorded selection after making the edit.
$ws.Range("C10:C39").Select()
